$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BurndownChart")

# Sprint 2 INSERT klant
# "artikelen bestellen" (row 11) and "inkooporders beheren" (row 12) are no
# longer marked as finished in wk21 (column E) -> clear the "x" marker.
$ws.Range("E11").ClearContents()
$ws.Range("E12").ClearContents()

# "klant toevoegen" (row 13) is the newly INSERTed task and is finished in
# wk21 with its 6 story points logged.
$ws.Range("E13").Value = 6

# "klant gegeven opzoeken" (row 20) also got 2 of its points done in wk21.
$ws.Range("E20").Value = 2

# Move the active selection to reflect where the edit happened.
$ws.Range("E21").Select()

# Recalculate so the "Nog doen" totals (row 24) pick up the new remaining
# points before the workbook is saved.
$excel.CalculateFullRebuild()

$wb.Save()
